# Script 1 - atualização automática de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g1.2")

$ws.Range("A14").Value = 2023
$ws.Range("B14").Value = 3.241657824791806
$ws.Range("C14").Value = 2.867008788862638
$ws.Range("D14").Value = 3.118144130554446
